# Update the "Estrategias_Optimas" data table with the newly computed rows.
# The original table had 10 data rows (A2:G11). The new table has 20 data
# rows (A2:G21): 3 brand-new low-value rows are inserted at the top, the
# previous 10 rows shift down and get new intermediate rows mixed in, and
# 3 new rows are appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, 1, 1, 50, 4, 23000, 0.04040085616235199),
    @(0, 2, 2, 100, 4, 45000, 0.09766834365325078),
    @(0, 3, 3, 50, 4, 66400, 0.1056810634474137),
    @(0, 4, 7, 150, 4, 69400, 0.2164772727272727),
    @(0, 4, 6, 150, 4, 63200, 0.2145833333333333),
    @(0, 4, 6, 100, 4, 170200, 0.1921882832080201),
    @(0, 4, 5, 100, 4, 151800, 0.1890533625730994),
    @(0, 4, 8, 150, 4, 59400, 0.1770021645021645),
    @(0, 5, 6, 150, 4, 113800, 0.30625),
    @(0, 5, 5, 150, 4, 90200, 0.250297619047619),
    @(0, 5, 2, 50, 4, 125800, 0.1605156679822188),
    @(1, 2, 2, 100, 4, 19800, 0.1806308049535604),
    @(1, 4, 2, 100, 4, 84600, 0.3441057394617766),
    @(1, 5, 2, 100, 4, 94800, 0.4059061681352703),
    @(1, 12, 4, 150, 4, 574200, 0.6767144892144892),
    @(1, 12, 4, 200, 4, 574200, 0.6767144892144892),
    @(1, 12, 4, 250, 4, 574200, 0.6767144892144892),
    @(1, 12, 5, 150, 4, 639600, 0.6724780701754386),
    @(1, 12, 5, 200, 4, 639600, 0.6724780701754386),
    @(1, 12, 5, 250, 4, 639600, 0.6724780701754386)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
